$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# Section 1: rows 39-42 in column E (new ChatGPT thread entry)
# Insertion order matches the shared-string table order of the target
# file: url (8), "Burnaby Tool" (9), prompt (10), "rank and rent..." (11)
# ------------------------------------------------------------------
$ws.Range("E41").Value = "https://chatgpt.com/c/68b2540b-2584-8321-b9c9-9bfaa8ac6050"
$ws.Range("E39").Value = "Burnaby Tool"
$ws.Range("E42").Value = "please create midjourney prompts for me to use to create logos for my rank and rent websites about duct cleaning they must include no text, be somewhat minimalistic and logo looking. and generally have a transparent bg or can be transitioned easily to a png through a converter tool (from the layout / background of the logo)"
$ws.Range("E40").Value = "rank and rent site image logos"
$ws.Range("E40").Font.Bold = $true

# ------------------------------------------------------------------
# Section 2: row 55 (city / state "dp1" table header row), left to right
# ------------------------------------------------------------------
$ws.Range("A55").Value = "ann arbor"
$ws.Range("A55").HorizontalAlignment = -4131

$ws.Range("B55").Value = "f"
$ws.Range("C55").Value = "flint"
$ws.Range("D55").Value = "f"
$ws.Range("F55").Value = "f"
$ws.Range("H55").Value = "f"
$ws.Range("I55").Value = "jackson"
$ws.Range("J55").Value = "f"
$ws.Range("K55").Value = "tuscaloosa"
$ws.Range("L55").Interior.ThemeColor = 1
$ws.Range("M55").Value = "peoria"
$ws.Range("N55").Interior.ThemeColor = 1
$ws.Range("O55").Value = "springfield"
$ws.Range("P55").Interior.ThemeColor = 1
$ws.Range("Q55").Value = "chattanooga"
$ws.Range("R55").Interior.ThemeColor = 1
$ws.Range("S55").Value = "charleston"
$ws.Range("T55").Interior.ThemeColor = 1
$ws.Range("U55").Value = "augusta"
$ws.Range("V55").Value = "f"
$ws.Range("W55").Value = "augusta"
$ws.Range("X55").Value = "f"
$ws.Range("Y55").Value = "augusta"
$ws.Range("Z55").Value = "f"
$ws.Range("AB55").Value = "f"
$ws.Range("AC55").Value = "flint"
$ws.Range("AD55").Value = "f"
$ws.Range("AE55").Value = "ft myers"
$ws.Range("AF55").Value = "f"
$ws.Range("AG55").Value = "cape coral"
$ws.Range("AH55").Value = "f"
$ws.Range("AI55").Value = "pembroke pines"
$ws.Range("AJ55").Value = "f"
$ws.Range("AK55").Value = "savannha"

# ------------------------------------------------------------------
# Section 3: row 56 (state-abbreviation row), left to right
# ------------------------------------------------------------------
$ws.Range("A56").Value = "mi"
$ws.Range("A56").HorizontalAlignment = -4131

$ws.Range("B56").Value = "f"
$ws.Range("C56").Value = "mi"
$ws.Range("D56").Value = "f"
$ws.Range("F56").Value = "f"
$ws.Range("H56").Value = "f"
$ws.Range("I56").Value = "wy"
$ws.Range("J56").Value = "f"
$ws.Range("K56").Value = "al"
$ws.Range("L56").Interior.ThemeColor = 1
$ws.Range("M56").Value = "il"
$ws.Range("N56").Interior.ThemeColor = 1
$ws.Range("O56").Value = "mo"
$ws.Range("P56").Interior.ThemeColor = 1
$ws.Range("Q56").Value = "tn"
$ws.Range("R56").Interior.ThemeColor = 1
$ws.Range("S56").Value = "sc"
$ws.Range("T56").Interior.ThemeColor = 1
$ws.Range("U56").Value = "ga"
$ws.Range("V56").Value = "f"
$ws.Range("W56").Value = "ga"
$ws.Range("X56").Value = "f"
$ws.Range("Y56").Value = "ga"
$ws.Range("Z56").Value = "f"
$ws.Range("AB56").Value = "f"
$ws.Range("AC56").Value = "mi"
$ws.Range("AD56").Value = "f"
$ws.Range("AE56").Value = "fl"
$ws.Range("AF56").Value = "f"
$ws.Range("AG56").Value = "fl"
$ws.Range("AH56").Value = "f"
$ws.Range("AI56").Value = "fl"
$ws.Range("AJ56").Value = "f"
$ws.Range("AK56").Value = "ga"

# ------------------------------------------------------------------
# Section 4: black-background ("f" spacer) cell fills, reusing the
# existing black (theme 1) fill wherever the engine allows it.
# ------------------------------------------------------------------
$blackCells = @(
    "B55","D55","F55","H55","J55","V55","X55","Z55","AB55","AD55","AF55","AH55","AJ55",
    "B56","D56","F56","H56","J56","V56","X56","Z56","AB56","AD56","AF56","AH56","AJ56"
)
foreach ($cellRef in $blackCells) {
    $ws.Range($cellRef).Interior.ThemeColor = 1
}
